$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Zeiterfassung")

# Insert a new row at row 9 (pushes the existing empty placeholder rows down by one,
# and expands the adjacent table "Tabelle1" by one row).
$ws.Rows.Item(9).Insert()

# Fill in the two new "Grobplanung" (rough planning) estimate rows.
$ws.Range("B9").Value = "PSP Draf, Ausserorderntliche Planung"
$ws.Range("C9").Value = "Tobias Lanz"
$ws.Range("D9").Value = 42259
$ws.Range("F9").Value = 4

$ws.Range("B10").Value = "Grob Planung - Lastenhef Reviews"
$ws.Range("C10").Value = "Tobias Lanz"
$ws.Range("D10").Value = 42260
$ws.Range("F10").Value = 4

# Expand the "Tabelle1" table by one row to account for the inserted row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B1:F101"))

$ws.Range("L35").Select()
